# Auto-generated Excel COM-interop script
# Applies numeric corrections to H:N columns across ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 12658.2
$ws.Range("I18").Value = 13611.111
$ws.Range("K18").Value = 13611.111
$ws.Range("M18").Value = -13327.111

$ws.Range("H53").Value = 403.5
$ws.Range("J53").Value = 288.16666
$ws.Range("L53").Value = 288.16666
$ws.Range("N53").Value = -1562.16666

$ws.Range("H97").Value = 1410.875
$ws.Range("I97").Value = 700
$ws.Range("K97").Value = 2100
$ws.Range("M97").Value = -1604

$ws.Range("H99").Value = 90924620
$ws.Range("J99").Value = 200000620
$ws.Range("L99").Value = 600001860
$ws.Range("N99").Value = -600004856

$ws.Range("H100").Value = 5871.85
$ws.Range("I100").Value = 5555.1177
$ws.Range("J100").Value = 7666.6665
$ws.Range("K100").Value = 5555.1177
$ws.Range("L100").Value = 7666.6665
$ws.Range("M100").Value = -5014.1177
$ws.Range("N100").Value = -8748.666499999999

$ws.Range("I111").Value = 357
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1071
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = 1996

$ws.Range("H125").Value = 1377.4667
$ws.Range("I125").Value = 986.4
$ws.Range("J125").Value = 2159.6
$ws.Range("K125").Value = 8877.6
$ws.Range("L125").Value = 19436.4
$ws.Range("M125").Value = -6417.6
$ws.Range("N125").Value = -24356.4

$ws.Range("H131").Value = 9923.529
$ws.Range("J131").Value = 75500
$ws.Range("L131").Value = 226500
$ws.Range("N131").Value = -236580

$ws.Range("H132").Value = 2291.3403
$ws.Range("I132").Value = 2368.0977
$ws.Range("K132").Value = 7104.293099999999
$ws.Range("M132").Value = -4574.293099999999

$ws.Range("H135").Value = 1037.4814
$ws.Range("I135").Value = 1190.8572
$ws.Range("K135").Value = 10717.7148
$ws.Range("M135").Value = -8182.7148

$ws.Range("H137").Value = 1518.4445
$ws.Range("I137").Value = 1402.8334
$ws.Range("K137").Value = 4208.5002
$ws.Range("M137").Value = -1658.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3667.3064
$ws.Range("I32").Value = 3078.125
$ws.Range("K32").Value = 3078.125
$ws.Range("M32").Value = -2791.125

$ws.Range("H97").Value = 1814
$ws.Range("I97").Value = 1814
$ws.Range("K97").Value = 1814
$ws.Range("M97").Value = -1318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 18043.846
$ws.Range("I99").Value = 18043.846
$ws.Range("K99").Value = 18043.846
$ws.Range("M99").Value = -16545.846

$ws.Range("H107").Value = 3560.9285
$ws.Range("I107").Value = 1442.75
$ws.Range("J107").Value = 6385.1665
$ws.Range("K107").Value = 1442.75
$ws.Range("L107").Value = 6385.1665
$ws.Range("M107").Value = 477.25
$ws.Range("N107").Value = -10225.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1000
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -887

$ws.Range("H31").Value = 9516.439
$ws.Range("I31").Value = 1681.4193
$ws.Range("K31").Value = 1681.4193
$ws.Range("M31").Value = -1386.4193

$ws.Range("H34").Value = 9516.439
$ws.Range("I34").Value = 1681.4193
$ws.Range("K34").Value = 1681.4193
$ws.Range("M34").Value = -1479.4193

$ws.Range("H50").Value = 14999.857
$ws.Range("J50").Value = 14999.857
$ws.Range("L50").Value = 14999.857
$ws.Range("N50").Value = -16249.857

$ws.Range("H51").Value = 15000
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -16472

$ws.Range("H59").Value = 19998.334
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290

$ws.Range("H60").Value = 15000
$ws.Range("J60").Value = 15000
$ws.Range("L60").Value = 15000
$ws.Range("N60").Value = -16022

$ws.Range("H61").Value = 15000
$ws.Range("J61").Value = 15000
$ws.Range("L61").Value = 15000
$ws.Range("N61").Value = -15696

$ws.Range("H88").Value = 19781
$ws.Range("J88").Value = 19781
$ws.Range("L88").Value = 19781
$ws.Range("N88").Value = -20593

$ws.Range("H91").Value = 19781
$ws.Range("J91").Value = 19781
$ws.Range("L91").Value = 19781
$ws.Range("N91").Value = -22589

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11000
$ws.Range("I70").Value = 11000
$ws.Range("K70").Value = 11000
$ws.Range("M70").Value = -10730

$ws.Range("H73").Value = 11000
$ws.Range("I73").Value = 11000
$ws.Range("K73").Value = 11000
$ws.Range("M73").Value = -10064

$ws.Range("H97").Value = 42895
$ws.Range("I97").Value = 82920.125
$ws.Range("J97").Value = 2869.875
$ws.Range("K97").Value = 82920.125
$ws.Range("L97").Value = 2869.875
$ws.Range("M97").Value = -82424.125
$ws.Range("N97").Value = -3861.875

$ws.Range("H102").Value = 2381.6
$ws.Range("I102").Value = 2408.8572
$ws.Range("K102").Value = 2408.8572
$ws.Range("M102").Value = -786.8571999999999

$ws.Range("H113").Value = 3913.077
$ws.Range("I113").Value = 3412.5715
$ws.Range("K113").Value = 3412.5715
$ws.Range("M113").Value = -1242.5715

$ws.Range("H126").Value = 2192.3333
$ws.Range("I126").Value = 2192.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6576.999899999999
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -4106.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1715.1666
$ws.Range("I22").Value = 2073.25
$ws.Range("K22").Value = 2073.25
$ws.Range("M22").Value = -1778.25

$ws.Range("H27").Value = 1715.1666
$ws.Range("I27").Value = 2073.25
$ws.Range("K27").Value = 2073.25
$ws.Range("M27").Value = -1966.25

$ws.Range("H40").Value = 9085.15
$ws.Range("I40").Value = 7911.3335
$ws.Range("K40").Value = 7911.3335
$ws.Range("M40").Value = -7775.3335

$ws.Range("H46").Value = 31786.072
$ws.Range("I46").Value = 106778
$ws.Range("K46").Value = 106778
$ws.Range("M46").Value = -106590

$ws.Range("H55").Value = 661.1667
$ws.Range("I55").Value = 778.3333
$ws.Range("J55").Value = 544
$ws.Range("K55").Value = 778.3333
$ws.Range("L55").Value = 544
$ws.Range("M55").Value = -605.3333
$ws.Range("N55").Value = -890

$ws.Range("H68").Value = 2249.75
$ws.Range("I68").Value = 1999.5
$ws.Range("K68").Value = 1999.5
$ws.Range("M68").Value = -1250.5

$ws.Range("H71").Value = 2249.75
$ws.Range("I71").Value = 1999.5
$ws.Range("K71").Value = 9997.5
$ws.Range("M71").Value = -6253.5

$ws.Range("H82").Value = 2669.875
$ws.Range("I82").Value = 2373.2
$ws.Range("J82").Value = 3164.3333
$ws.Range("K82").Value = 2373.2
$ws.Range("L82").Value = 3164.3333
$ws.Range("M82").Value = -2012.2
$ws.Range("N82").Value = -3886.3333

$ws.Range("H85").Value = 2669.875
$ws.Range("I85").Value = 2373.2
$ws.Range("J85").Value = 3164.3333
$ws.Range("K85").Value = 2373.2
$ws.Range("L85").Value = 3164.3333
$ws.Range("M85").Value = -1125.2
$ws.Range("N85").Value = -5660.3333

$ws.Range("H93").Value = 34076.363
$ws.Range("I93").Value = 3578.75
$ws.Range("K93").Value = 3578.75
$ws.Range("M93").Value = -2330.75

$ws.Range("H132").Value = 4133.926
$ws.Range("I132").Value = 3601.524
$ws.Range("J132").Value = 5997.3335
$ws.Range("K132").Value = 10804.572
$ws.Range("L132").Value = 17992.0005
$ws.Range("M132").Value = -8274.572
$ws.Range("N132").Value = -23052.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3219.8
$ws.Range("I81").Value = 3274.75
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 6549.5
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -5488.5
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 3219.8
$ws.Range("I84").Value = 3274.75
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 32747.5
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -27443.5
$ws.Range("N84").Value = -40608

$ws.Range("H107").Value = 20006734
$ws.Range("I107").Value = 10166.75
$ws.Range("J107").Value = 38465104
$ws.Range("K107").Value = 30500.25
$ws.Range("L107").Value = 115395312
$ws.Range("M107").Value = -28580.25
$ws.Range("N107").Value = -115399152

$ws.Range("H136").Value = 8759.6
$ws.Range("I136").Value = 8621.777
$ws.Range("K136").Value = 25865.331
$ws.Range("M136").Value = -23315.331
